$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 18:05"

# Country data rows: some countries swapped rank (label + stats move together)
# and/or received refreshed case counts. Each entry is the full final row
# (columns A-H) for that row number.
$rowsData = @(
    @{ Row = 4; Values = @("Estados Unidos", 1671658, 4830, 447056, 1125827, 0, 92, 98775) },
    @{ Row = 5; Values = @("Brasil", 352163, 4765, 142587, 187362, 0, 201, 22214) },
    @{ Row = 11; Values = @("Alemania", 180084, 98, 160300, 11413, 0, 5, 8371) },
    @{ Row = 13; Values = @("India", 137608, 6185, 54865, 78739, 0, 136, 4004) },
    @{ Row = 14; Values = @("Iran", 135701, 2180, 105801, 22483, 0, 58, 7417) },
    @{ Row = 36; Values = @("Polonia", 21326, 395, 9194, 11136, 0, 3, 996) },
    @{ Row = 37; Values = @("Kuwait", 21302, 838, 6117, 15029, 0, 8, 156) },
    @{ Row = 54; Values = @("Chequia", 8932, 42, 6063, 2554, 0, 1, 315) },
    @{ Row = 71; Values = @("Luxemburgo", 3992, 2, 3767, 115, 0, 1, 110) },
    @{ Row = 81; Values = @("Grecia", 2878, 2, 1374, 1333, 0, 0, 171) },
    @{ Row = 143; Values = @("Mauricio", 334, 2, 322, 2, 0, 0, 10) },
    @{ Row = 187; Values = @("Botsuana", 35, 5, 19, 15, 0, 0, 1) },
    @{ Row = 188; Values = @("Guam", 32, 0, 0, 31, 0, 0, 1) },
    @{ Row = 198; Values = @("Nueva Caledonia", 18, 0, 18, 0, 0, 0, 0) },
    @{ Row = 199; Values = @("Belice", 18, 0, 16, 0, 0, 0, 2) },
    @{ Row = 200; Values = @("Santa Lucia", 18, 0, 18, 0, 0, 0, 0) },
    @{ Row = 209; Values = @("Seychelles", 11, 0, 11, 0, 0, 0, 0) },
    @{ Row = 210; Values = @("Groenlandia", 11, 0, 11, 0, 0, 0, 0) },
    @{ Row = 211; Values = @("Montserrat", 11, 0, 10, 0, 0, 0, 1) },
    @{ Row = 214; Values = @("Bonaire, San Eustaquio y Saba", 6, 0, 6, 0, 0, 0, 0) },
    @{ Row = 216; Values = @("Sahara Occidental", 6, 0, 6, 0, 0, 0, 0) }
)

foreach ($item in $rowsData) {
    $r = $item.Row
    $vals = $item.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
